$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 735 ("「アルマジロ」" entry) entirely; this shifts all
# subsequent rows up by one, matching the renumbering seen in the diff.
$ws.Rows(735).Delete()
